$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns for data rows are treated as text so that
# numeric-looking strings (e.g. "1.000") are not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '31.304.51'
$ws.Range("E2").Value = '  +2.45%  '

# Row 3: Ethereum
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.967.78'
$ws.Range("E3").Value = '  +2.83%  '

# Row 4: TetherUSD
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5: BNB
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '247.85'
$ws.Range("E5").Value = '  +1.36%  '

# Row 6: USDC
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.06%  '

# Row 7: XRP
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.4898'
$ws.Range("E7").Value = '  +0.99%  '

# Row 8: OKB
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '44.77'
$ws.Range("E8").Value = '  +0.72%  '

# Row 9: Cardano
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2986'
$ws.Range("E9").Value = '  +3.27%  '

# Row 10: Dogecoin
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.06857'
$ws.Range("E10").Value = '  +0.75%  '

# Row 11: Solana
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '19.30'
$ws.Range("E11").Value = '  -0.30%  '

# Row 12: Litecoin
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '107.53'
$ws.Range("E12").Value = '  -3.23%  '

# Row 13: WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.955.01'
$ws.Range("E13").Value = '  +2.09%  '

# Row 14: TRON
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.07775'
$ws.Range("E14").Value = '  +2.71%  '

# Row 15: Polkadot
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.470'
$ws.Range("E15").Value = '  +1.65%  '

# Row 16: Polygon
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.7177'
$ws.Range("E16").Value = '  +6.98%  '

# Row 17: BitcoinCash
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '288.52'
$ws.Range("E17").Value = '  -2.58%  '

# Row 18: WrappedBTC
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '31.321.56'
$ws.Range("E18").Value = '  +2.50%  '

# Row 19: Avalanche
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '13.33'
$ws.Range("E19").Value = '  +2.44%  '

# Row 20: ShibaInu
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007768'
$ws.Range("E20").Value = '  +2.31%  '

# Row 21: Uniswap
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '5.644'
$ws.Range("E21").Value = '  +2.03%  '

# Row 22: WrappedliquidstakedEther2.0
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.200.80'
$ws.Range("E22").Value = '  +1.74%  '

# Row 23: Dai
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.03%  '

# Row 24: BinanceUSD
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").Value = '  -0.02%  '

# Row 25: Chainlink
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '6.648'
$ws.Range("E25").Value = '  +3.08%  '

# Row 26: Cosmos
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '10.05'
$ws.Range("E26").Value = '  +6.12%  '

# Row 27: Monero
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '169.20'
$ws.Range("E27").Value = '  +1.96%  '

# Row 28: EthereumClassic
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '20.07'
$ws.Range("E28").Value = '  -1.17%  '

# Row 29: LidoDAOToken
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.193'
$ws.Range("E29").Value = '  +5.52%  '

# Row 30: Stellar
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.1070'
$ws.Range("E30").Value = '  +0.26%  '

# Row 31: Toncoin
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.449'
$ws.Range("E31").Value = '  +1.04%  '

# Row 32: Filecoin
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.884'
$ws.Range("E32").Value = '  +20.44%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.521'
$ws.Range("E33").Value = '  +8.97%  '

# Row 34: Hedera
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05076'
$ws.Range("E34").Value = '  +1.79%  '

# Row 35: ImmutableX
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7726'
$ws.Range("E35").Value = '  +5.13%  '

# Row 36: ARBITRUM
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.174'
$ws.Range("E36").Value = '  +2.95%  '

# Row 37: VeChain
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02057'
$ws.Range("E37").Value = '  +1.02%  '

# Row 38: HuobiToken
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.730'
$ws.Range("E38").Value = '  +0.46%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.715'
$ws.Range("E39").Value = '  +1.13%  '

# Row 40: RenderToken
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.146'
$ws.Range("E40").Value = '  +6.26%  '

# Row 41: FraxShare
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.401'
$ws.Range("E41").Value = '  +10.10%  '

# Row 42: Aave
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '73.99'
$ws.Range("E42").Value = '  +6.34%  '

# Row 43: Quant
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '110.09'
$ws.Range("E43").Value = '  +0.86%  '

# Row 44: TrustWalletToken
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.8864'
$ws.Range("E44").Value = '  +2.19%  '

# Row 45: TheSandbox
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4488'
$ws.Range("E45").Value = '  +0.91%  '

# Row 46: PaxDollar
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9995'
$ws.Range("E46").Value = '  -0.03%  '

# Row 47: Aptos
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.534'
$ws.Range("E47").Value = '  +4.50%  '

# Row 48: Maker
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '997.14'
$ws.Range("E48").Value = '  +18.16%  '

# Row 49: EnergySwap
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.440'
$ws.Range("E49").Value = '  +2.02%  '

# Row 50: Algorand
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1272'
$ws.Range("E50").Value = '  +3.78%  '

# Row 51: Elrond
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '36.14'
$ws.Range("E51").Value = '  +4.14%  '
